# Update "想去人数" (F column) figures for both the "展览" sheet and the
# "全部类型" sheet, which mirror the same rows of data.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 14795
    3  = 18255
    15 = 190
    17 = 1381
    22 = 7559
    27 = 14
    28 = 5915
    34 = 5241
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    2  = 14795
    3  = 18255
    15 = 190
    17 = 1381
    23 = 7559
    28 = 14
    30 = 5915
    36 = 5241
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
